# Maximum Capacity Factor.xlsx update
# - "About" sheet: bump the last-updated date in C1
# - "MCF" sheet: set all (non-100%) capacity-factor inputs to 100% (1)
#   The dependent formula cells (B19:B22, B24:B25) recalc automatically.
# - Leave the active selection on the MCF sheet at B17, matching the
#   cell last touched by the edit.

$wb = $excel.ActiveWorkbook

# --- About sheet: update the date stamp in C1 ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45392

# --- MCF sheet: bump capacity factors to 100% ---
$mcf = $wb.Worksheets.Item("MCF")

$cellsToMax = @("B2", "B3", "B4", "B6", "B10", "B11", "B12", "B13", "B14", "B16", "B17", "B18")
foreach ($addr in $cellsToMax) {
    $mcf.Range($addr).Value = 1
}

# Leave the selection where the user last edited
$mcf.Activate()
$mcf.Range("B17").Select()
